$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M13").Value = 2180.66

# Sheet "VENTA MENSUAL"
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F13").Value = 2291.68
$wsVentaMensual.Range("F22").Value = 16273.63

# Sheet "CUMPLIMIENTO MENSUAL"
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimientoMensual.Range("D16").Value = 5692.99
$wsCumplimientoMensual.Range("E16").Value = 23839.45
$wsCumplimientoMensual.Range("F16").Value = 0.1927707294080679

$wsCumplimientoMensual.Range("D19").Value = 16273.63
$wsCumplimientoMensual.Range("E19").Value = 34113.56762291768
$wsCumplimientoMensual.Range("F19").Value = 0.3229715238737199
